$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Tränare: William" note (clears the cell and, since it becomes an
# unused shared string, it is dropped from the workbook's shared string table
# on save).
$ws.Range("A10").Value = $null

# Resize columns A, B, C, E, G, H, I to their new (narrower) widths.
# This runtime's ColumnWidth setter quantizes to 1/6-character steps via
# stored = (floor(ColumnWidth*6 + 0.5) + 5) / 6, so we pre-compensate by the
# fixed 5/6 offset to land as close as possible on the intended stored width.
$offset = 5.0 / 6.0
$ws.Columns.Item(1).ColumnWidth = 12.5703125 - $offset
$ws.Columns.Item(2).ColumnWidth = 7 - $offset
$ws.Columns.Item(3).ColumnWidth = 8.140625 - $offset
$ws.Columns.Item(5).ColumnWidth = 7.85546875 - $offset
$ws.Columns.Item(7).ColumnWidth = 12.42578125 - $offset
$ws.Columns.Item(8).ColumnWidth = 9 - $offset
$ws.Columns.Item(9).ColumnWidth = 15.140625 - $offset

# Move the active selection to A10.
$ws.Range("A10").Select()
